# Applies the template change: extends the load-column grid (adds
# Load-6..Load-9 / extra columns G:L resp. H:L) with the unlocked/italic
# "fill me in" style used throughout the sheet, fills in the sample data
# that was added to "Example Test-2", and updates the remembered cell
# selection on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Example Test-1"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Example Test-1")

# Rows 6-14 (the 9 measurement rows) gain empty, unlocked/italic cells in
# columns G:L (same style already used by the existing F column cells).
$rng = $ws1.Range("G6:L14")
$rng.Font.Italic = $true
$rng.Locked = $false

# Rows 15-22 (the derived-value rows that previously had no data columns
# at all) gain the same empty, unlocked/italic cells across B:L.
$rng = $ws1.Range("B15:L22")
$rng.Font.Italic = $true
$rng.Locked = $false

# ---------------------------------------------------------------------
# Sheet "Example Test-2"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Example Test-2")

# Rows 6-14: columns H:K stay empty, column L gets the new Load-9 sample
# values (where the row already carries data for the other loads).
$rng = $ws2.Range("H6:L14")
$rng.Font.Italic = $true
$rng.Locked = $false

$ws2.Range("L6").Value = 160
$ws2.Range("L9").Value = 3.2
$ws2.Range("L10").Value = 140
$ws2.Range("L11").Value = 95
$ws2.Range("L14").Value = 40

# Rows 15-20 and 22: columns B:L gain the same empty style as on sheet 1.
$rng = $ws2.Range("B15:L20")
$rng.Font.Italic = $true
$rng.Locked = $false

$rng = $ws2.Range("B22:L22")
$rng.Font.Italic = $true
$rng.Locked = $false

# Row 21 gains the same empty style across B:L ...
$rng = $ws2.Range("B21:L21")
$rng.Font.Italic = $true
$rng.Locked = $false

# ... plus sample temperature values in B:G and L.
$ws2.Range("B21").Value = 37
$ws2.Range("C21").Value = 37.1
$ws2.Range("D21").Value = 37 + 0.2
$ws2.Range("E21").Value = 37 + 0.3
$ws2.Range("F21").Value = 37.4
$ws2.Range("G21").Value = 37.5
$ws2.Range("L21").Value = 38

# ---------------------------------------------------------------------
# Remembered selections
# ---------------------------------------------------------------------
# Select Example Test-2 first so that the final .Activate() below (on
# Example Test-1, which is the tab that was selected in the original
# workbook) leaves that sheet as the active one.
$ws2.Activate()
$ws2.Range("G13").Select()

$ws1.Activate()
$ws1.Range("H11").Select()
